$d = $word.ActiveDocument
$wXml = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Table 2, row 3: fill in the "Search Dictionary Function" test entry ---
$tbl = $d.Tables.Item(2)

$tbl.Cell(3, 1).Range.Text = "19/6/25"
$tbl.Cell(3, 2).Range.Text = "Search Dictionary Function"
$tbl.Cell(3, 3).Range.Text = "Expected"
$tbl.Cell(3, 4).Range.Text = "Task_dictionary, “T1”"
$tbl.Cell(3, 5).Range.Text = "The function returns the details of task “T1”"

# Column 6 ("Actual Result") holds two runs split by a rendered page break.
$xml6 = "<w:p $wXml><w:r><w:t xml:space=`"preserve`">The function returns the details, but no </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>formatting has been applied</w:t></w:r></w:p>"
$tbl.Cell(3, 6).Range.InsertXML($xml6)
$tbl = $d.Tables.Item(2)
$tbl.Cell(3, 6).Range.Paragraphs.Item(1).Range.Delete()

# Column 7 ("How it was fixed") holds two runs, each preceded by a rendered page break.
$tbl = $d.Tables.Item(2)
$xml7 = "<w:p $wXml><w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">A function should be made for single dictionary formatting where a </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>nested system is not needed</w:t></w:r></w:p>"
$tbl.Cell(3, 7).Range.InsertXML($xml7)
$tbl = $d.Tables.Item(2)
$tbl.Cell(3, 7).Range.Paragraphs.Item(1).Range.Delete()

# --- Table 4, row 2, cell 1: mark the "Date" header with a rendered page break ---
$tbl4 = $d.Tables.Item(4)
$xmlDate = "<w:p $wXml><w:r><w:lastRenderedPageBreak/><w:t>Date</w:t></w:r></w:p>"
$tbl4.Cell(2, 1).Range.InsertXML($xmlDate)
$tbl4 = $d.Tables.Item(4)
$tbl4.Cell(2, 1).Range.Paragraphs.Item(1).Range.Delete()
